$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.996.23"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "2.743.28"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'572.56"
$ws.Range("E5").Value = "  -1.08%  "
$ws.Range("D6").Value = "'160.68"
$ws.Range("E6").Value = "  +1.21%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'0.596"
$ws.Range("E8").Value = "  -1.90%  "
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("E10").Value = "  +4.95%  "
$ws.Range("D11").Value = "'5.77"
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("D12").Value = "'0.384"
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("D13").Value = "3.230.19"
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("D14").Value = "'27.01"
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("D15").Value = "63.843.56"
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("D16").Value = "'0.0000150"
$ws.Range("E16").Value = "  -1.63%  "
$ws.Range("D17").Value = "2.746.40"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").Value = "'12.23"
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("D19").Value = "'4.80"
$ws.Range("E19").Value = "  -1.84%  "
$ws.Range("D20").Value = "'354.20"
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("D21").Value = "'6.63"
$ws.Range("E21").Value = "  -3.41%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("E23").Value = "  -5.03%  "
$ws.Range("D24").Value = "'64.41"
$ws.Range("E24").Value = "  -2.25%  "
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("D28").Value = "0.0₃0914"
$ws.Range("E28").Value = "  -1.84%  "
$ws.Range("E29").Value = "  +0.99%  "
$ws.Range("E30").Value = "  +2.80%  "
$ws.Range("D31").Value = "'1.33"
$ws.Range("E31").Value = "  +7.71%  "
$ws.Range("D32").Value = "'163.97"
$ws.Range("E32").Value = "  -1.86%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'20.15"
$ws.Range("E33").Value = "  -1.06%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "'4.91"
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("D37").Value = "'1.80"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "'0.989"
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("D39").Value = "'349.58"
$ws.Range("D40").Value = "'6.39"
$ws.Range("E40").Value = "  +2.54%  "
$ws.Range("D41").Value = "'4.10"
$ws.Range("E41").Value = "  -1.43%  "
$ws.Range("E42").Value = "  -1.75%  "
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("D44").Value = "'21.18"
$ws.Range("E45").Value = "  -1.58%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'135.25"
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.625"
$ws.Range("E47").Value = "  -1.83%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0251"
$ws.Range("E48").Value = "  -2.88%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "'0.100"
$ws.Range("E49").Value = "  -1.24%  "
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "'11.06"
$ws.Range("E51").Value = "  +0.02%  "
